# Collapse the split "<id>" / "p049r_aN" / "</id>" runs into a single
# run reading "<id>p049r_N</id>" for each of the four occurrences
# (p049r_a1..a4 -> p049r_1..4) in the document.

$d = $word.ActiveDocument

for ($i = 1; $i -le 4; $i++) {
    $old = "<id>p049r_a$i</id>"
    $new = "<id>p049r_$i</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
